$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap values in C2 and C3 (data3.xlsx diff: C2 1->2, C3 2->1)
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 1
